$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Case/Study query (row 2, CasesTab): the "Cohort" coalesce line (and its
# matching "co.cohort_description" output) is no longer part of the query.
$caseQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nMATCH (c)<--(diag:diagnosis)`nMATCH (samp:sample)-->(c) `n WHERE samp.specific_sample_pathology IN [`"Oligodendroglioma`"]  `nOPTIONAL MATCH (co:cohort)<-[*]-(c)`n  WITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"
# Sample query (row 3, SamplesTab) -- unchanged text, rewritten so the shared
# string is re-inserted after the Case/Study query is updated.
$sampleQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) `n WHERE samp.specific_sample_pathology IN [`"Oligodendroglioma`"]  `nWITH DISTINCT samp AS samp, c, demo, diag`nRETURN  coalesce(samp.sample_id, '') AS ``Sample ID``, `n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(demo.breed,'') AS Breed , `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(samp.sample_site, '') AS ``Sample Site``,`n        coalesce(samp.summarized_sample_type, '') AS ``Sample Type``,`n        coalesce(samp.specific_sample_pathology, '') AS ``Pathology/Morphology``,`n        coalesce(samp.tumor_grade, '') AS ``Tumor Grade``,`n        coalesce(samp.sample_chronology, '') AS ``Sample Chronology``,`n        coalesce(samp.percentage_tumor, '') AS ``Percentage Tumor``,`n        coalesce(samp.necropsy_sample, '') AS ``Necropsy Sample``,`n        coalesce(samp.sample_preservation, '') AS ``Sample Preservation``"
# File query (row 4, FilesTab) -- unchanged text, rewritten for the same reason.
$fileQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n MATCH (samp:sample)-->(c) `n WHERE samp.specific_sample_pathology IN [`"Oligodendroglioma`"]  `nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B2").Value = $caseQuery
$ws.Range("B3").Value = $sampleQuery
$ws.Range("B4").Value = $fileQuery

# The author's selection/scroll position moved from B4 back up to B2.
[void]$ws.Range("B2").Select()
